# Auto-generated edit script: updates market-price derived cells
# across the per-job Sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) to match
# the latest scheduled-runner price snapshot.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 139.66667  # H38
$ws.Cells.Item(38, 9).Value = 139.66667  # I38
$ws.Cells.Item(38, 11).Value = 419.00001  # K38
$ws.Cells.Item(38, 13).Value = -47.00001000000003  # M38
$ws.Cells.Item(43, 8).Value = 6000  # H43
$ws.Cells.Item(43, 9).Value = 0  # I43
$ws.Cells.Item(43, 11).Value = 0  # K43
$ws.Cells.Item(43, 13).ClearContents()  # M43
$ws.Cells.Item(98, 8).Value = 766.25  # H98
$ws.Cells.Item(98, 9).Value = 446.2857  # I98
$ws.Cells.Item(98, 11).Value = 446.2857  # K98
$ws.Cells.Item(98, 13).Value = 1051.7143  # M98
$ws.Cells.Item(113, 8).Value = 17125.75  # H113
$ws.Cells.Item(113, 9).Value = 12352  # I113
$ws.Cells.Item(113, 11).Value = 12352  # K113
$ws.Cells.Item(113, 13).Value = -9098  # M113
$ws.Cells.Item(122, 8).Value = 766.25  # H122
$ws.Cells.Item(122, 9).Value = 446.2857  # I122
$ws.Cells.Item(122, 11).Value = 1338.8571  # K122
$ws.Cells.Item(122, 13).Value = 1111.1429  # M122
$ws.Cells.Item(140, 8).Value = 69420  # H140
$ws.Cells.Item(140, 10).Value = 69420  # J140
$ws.Cells.Item(140, 12).Value = 69420  # L140
$ws.Cells.Item(140, 14).Value = -79780  # N140

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 34998  # H44
$ws.Cells.Item(44, 10).Value = 34998  # J44
$ws.Cells.Item(44, 12).Value = 34998  # L44
$ws.Cells.Item(44, 14).Value = -35974  # N44
$ws.Cells.Item(55, 8).Value = 29000  # H55
$ws.Cells.Item(55, 10).Value = 29000  # J55
$ws.Cells.Item(55, 12).Value = 29000  # L55
$ws.Cells.Item(55, 14).Value = -29630  # N55
$ws.Cells.Item(74, 8).Value = 1323.6666  # H74
$ws.Cells.Item(74, 9).Value = 1323.6666  # I74
$ws.Cells.Item(74, 11).Value = 1323.6666  # K74
$ws.Cells.Item(74, 13).Value = -449.6666  # M74
$ws.Cells.Item(77, 8).Value = 1323.6666  # H77
$ws.Cells.Item(77, 9).Value = 1323.6666  # I77
$ws.Cells.Item(77, 11).Value = 6618.333000000001  # K77
$ws.Cells.Item(77, 13).Value = -2250.333000000001  # M77
$ws.Cells.Item(110, 8).Value = 979.6  # H110
$ws.Cells.Item(110, 9).Value = 974.5  # I110
$ws.Cells.Item(110, 11).Value = 974.5  # K110
$ws.Cells.Item(110, 13).Value = 1070.5  # M110
$ws.Cells.Item(122, 8).Value = 2888.375  # H122
$ws.Cells.Item(122, 9).Value = 2888.375  # I122
$ws.Cells.Item(122, 11).Value = 8665.125  # K122
$ws.Cells.Item(122, 13).Value = -6215.125  # M122

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 14947.5  # H35
$ws.Cells.Item(35, 10).Value = 14947.5  # J35
$ws.Cells.Item(35, 12).Value = 14947.5  # L35
$ws.Cells.Item(35, 14).Value = -15567.5  # N35
$ws.Cells.Item(82, 8).Value = 29149.6  # H82
$ws.Cells.Item(82, 9).Value = 16892  # I82
$ws.Cells.Item(82, 10).Value = 39875  # J82
$ws.Cells.Item(82, 11).Value = 16892  # K82
$ws.Cells.Item(82, 12).Value = 39875  # L82
$ws.Cells.Item(82, 13).Value = -16509  # M82
$ws.Cells.Item(82, 14).Value = -40641  # N82
$ws.Cells.Item(85, 8).Value = 29149.6  # H85
$ws.Cells.Item(85, 9).Value = 16892  # I85
$ws.Cells.Item(85, 10).Value = 39875  # J85
$ws.Cells.Item(85, 11).Value = 16892  # K85
$ws.Cells.Item(85, 12).Value = 39875  # L85
$ws.Cells.Item(85, 13).Value = -15566  # M85
$ws.Cells.Item(85, 14).Value = -42527  # N85
$ws.Cells.Item(86, 8).Value = 5626.6665  # H86
$ws.Cells.Item(86, 9).Value = 4925.875  # I86
$ws.Cells.Item(86, 11).Value = 4925.875  # K86
$ws.Cells.Item(86, 13).Value = -3802.875  # M86
$ws.Cells.Item(89, 8).Value = 5626.6665  # H89
$ws.Cells.Item(89, 9).Value = 4925.875  # I89
$ws.Cells.Item(89, 11).Value = 24629.375  # K89
$ws.Cells.Item(89, 13).Value = -19013.375  # M89
$ws.Cells.Item(105, 8).Value = 2727  # H105
$ws.Cells.Item(105, 9).Value = 2600.6365  # I105
$ws.Cells.Item(105, 11).Value = 2600.6365  # K105
$ws.Cells.Item(105, 13).Value = -853.6365000000001  # M105
$ws.Cells.Item(107, 8).Value = 1022.625  # H107
$ws.Cells.Item(107, 9).Value = 797.46155  # I107
$ws.Cells.Item(107, 11).Value = 797.46155  # K107
$ws.Cells.Item(107, 13).Value = 1122.53845  # M107

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 908.2857  # H16
$ws.Cells.Item(16, 9).Value = 901.3333  # I16
$ws.Cells.Item(16, 11).Value = 901.3333  # K16
$ws.Cells.Item(16, 13).Value = -614.3333  # M16
$ws.Cells.Item(41, 8).Value = 19998.75  # H41
$ws.Cells.Item(41, 10).Value = 19998.75  # J41
$ws.Cells.Item(41, 12).Value = 19998.75  # L41
$ws.Cells.Item(41, 14).Value = -20854.75  # N41
$ws.Cells.Item(50, 8).Value = 28509.125  # H50
$ws.Cells.Item(50, 10).Value = 29570  # J50
$ws.Cells.Item(50, 12).Value = 29570  # L50
$ws.Cells.Item(50, 14).Value = -30820  # N50
$ws.Cells.Item(51, 8).Value = 21772.5  # H51
$ws.Cells.Item(59, 8).Value = 34330.555  # H59
$ws.Cells.Item(59, 10).Value = 34621.875  # J59
$ws.Cells.Item(59, 12).Value = 34621.875  # L59
$ws.Cells.Item(59, 14).Value = -36911.875  # N59
$ws.Cells.Item(60, 8).Value = 22868.285  # H60
$ws.Cells.Item(60, 10).Value = 24997.5  # J60
$ws.Cells.Item(60, 12).Value = 24997.5  # L60
$ws.Cells.Item(60, 14).Value = -26019.5  # N60
$ws.Cells.Item(61, 8).Value = 21772.5  # H61
$ws.Cells.Item(107, 8).Value = 1244.8334  # H107
$ws.Cells.Item(107, 9).Value = 999  # I107
$ws.Cells.Item(107, 10).Value = 1294  # J107
$ws.Cells.Item(107, 11).Value = 999  # K107
$ws.Cells.Item(107, 12).Value = 1294  # L107
$ws.Cells.Item(107, 13).Value = 921  # M107
$ws.Cells.Item(107, 14).Value = -5134  # N107
$ws.Cells.Item(113, 8).Value = 908.2857  # H113
$ws.Cells.Item(113, 9).Value = 901.3333  # I113
$ws.Cells.Item(113, 11).Value = 901.3333  # K113
$ws.Cells.Item(113, 13).Value = 1268.6667  # M113
$ws.Cells.Item(122, 8).Value = 0  # H122
$ws.Cells.Item(122, 9).Value = 0  # I122
$ws.Cells.Item(122, 11).Value = 0  # K122
$ws.Cells.Item(122, 13).ClearContents()  # M122

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2166  # H34
$ws.Cells.Item(34, 10).Value = 4500  # J34
$ws.Cells.Item(34, 12).Value = 13500  # L34
$ws.Cells.Item(34, 14).Value = -13668  # N34
$ws.Cells.Item(39, 8).Value = 2500  # H39
$ws.Cells.Item(39, 10).Value = 0  # J39
$ws.Cells.Item(39, 12).Value = 0  # L39
$ws.Cells.Item(39, 14).ClearContents()  # N39
$ws.Cells.Item(55, 8).Value = 2849  # H55
$ws.Cells.Item(55, 9).Value = 1000  # I55
$ws.Cells.Item(55, 10).Value = 3054.4443  # J55
$ws.Cells.Item(55, 11).Value = 3000  # K55
$ws.Cells.Item(55, 12).Value = 9163.332900000001  # L55
$ws.Cells.Item(55, 13).Value = -2823  # M55
$ws.Cells.Item(55, 14).Value = -9517.332900000001  # N55

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 15165  # H43
$ws.Cells.Item(43, 10).Value = 15165  # J43
$ws.Cells.Item(43, 12).Value = 15165  # L43
$ws.Cells.Item(43, 14).Value = -15467  # N43
$ws.Cells.Item(70, 8).Value = 33338870  # H70
$ws.Cells.Item(70, 9).Value = 55559200  # I70
$ws.Cells.Item(70, 11).Value = 55559200  # K70
$ws.Cells.Item(70, 13).Value = -55558930  # M70
$ws.Cells.Item(73, 8).Value = 33338870  # H73
$ws.Cells.Item(73, 9).Value = 55559200  # I73
$ws.Cells.Item(73, 11).Value = 55559200  # K73
$ws.Cells.Item(73, 13).Value = -55558264  # M73
$ws.Cells.Item(102, 8).Value = 2961.6667  # H102
$ws.Cells.Item(102, 9).Value = 2998.5  # I102
$ws.Cells.Item(102, 10).Value = 2888  # J102
$ws.Cells.Item(102, 11).Value = 2998.5  # K102
$ws.Cells.Item(102, 12).Value = 2888  # L102
$ws.Cells.Item(102, 13).Value = -1376.5  # M102
$ws.Cells.Item(102, 14).Value = -6132  # N102
$ws.Cells.Item(109, 8).Value = 0  # H109
$ws.Cells.Item(109, 10).Value = 0  # J109
$ws.Cells.Item(109, 12).Value = 0  # L109
$ws.Cells.Item(109, 14).ClearContents()  # N109
$ws.Cells.Item(113, 8).Value = 2152.625  # H113
$ws.Cells.Item(113, 9).Value = 2144.4285  # I113
$ws.Cells.Item(113, 11).Value = 2144.4285  # K113
$ws.Cells.Item(113, 13).Value = 25.57150000000001  # M113
$ws.Cells.Item(122, 8).Value = 2855.7778  # H122
$ws.Cells.Item(122, 9).Value = 1615.3334  # I122
$ws.Cells.Item(122, 11).Value = 4846.0002  # K122
$ws.Cells.Item(122, 13).Value = -2396.0002  # M122

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 52877.8  # H7
$ws.Cells.Item(7, 9).Value = 52877.8  # I7
$ws.Cells.Item(7, 11).Value = 52877.8  # K7
$ws.Cells.Item(7, 13).Value = -52765.8  # M7
$ws.Cells.Item(40, 8).Value = 3000.5  # H40
$ws.Cells.Item(40, 9).Value = 3000.5  # I40
$ws.Cells.Item(40, 11).Value = 3000.5  # K40
$ws.Cells.Item(40, 13).Value = -2864.5  # M40
$ws.Cells.Item(55, 8).Value = 257.7  # H55
$ws.Cells.Item(55, 9).Value = 208  # I55
$ws.Cells.Item(55, 11).Value = 208  # K55
$ws.Cells.Item(55, 13).Value = -35  # M55
$ws.Cells.Item(122, 8).Value = 3463.2307  # H122
$ws.Cells.Item(122, 9).Value = 3151.8  # I122
$ws.Cells.Item(122, 10).Value = 4501.3335  # J122
$ws.Cells.Item(122, 11).Value = 9455.400000000001  # K122
$ws.Cells.Item(122, 12).Value = 13504.0005  # L122
$ws.Cells.Item(122, 13).Value = -7005.400000000001  # M122
$ws.Cells.Item(122, 14).Value = -18404.0005  # N122
$ws.Cells.Item(126, 8).Value = 52877.8  # H126
$ws.Cells.Item(126, 9).Value = 52877.8  # I126
$ws.Cells.Item(126, 11).Value = 158633.4  # K126
$ws.Cells.Item(126, 13).Value = -156163.4  # M126

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 29974.75  # H54
$ws.Cells.Item(54, 9).Value = 29899  # I54
$ws.Cells.Item(54, 10).Value = 30000  # J54
$ws.Cells.Item(54, 11).Value = 29899  # K54
$ws.Cells.Item(54, 12).Value = 30000  # L54
$ws.Cells.Item(54, 13).Value = -29379  # M54
$ws.Cells.Item(54, 14).Value = -31040  # N54
$ws.Cells.Item(113, 8).Value = 595.44446  # H113
$ws.Cells.Item(113, 9).Value = 506.2857  # I113
$ws.Cells.Item(113, 10).Value = 907.5  # J113
$ws.Cells.Item(113, 11).Value = 1518.8571  # K113
$ws.Cells.Item(113, 12).Value = 2722.5  # L113
$ws.Cells.Item(113, 13).Value = 651.1428999999998  # M113
$ws.Cells.Item(113, 14).Value = -7062.5  # N113
$ws.Cells.Item(122, 8).Value = 4131.1665  # H122
$ws.Cells.Item(122, 9).Value = 4131.1665  # I122
$ws.Cells.Item(122, 11).Value = 12393.4995  # K122
$ws.Cells.Item(122, 13).Value = -9943.499500000002  # M122
$ws.Cells.Item(126, 8).Value = 1693.75  # H126
$ws.Cells.Item(126, 10).Value = 1999  # J126
$ws.Cells.Item(126, 12).Value = 5997  # L126
$ws.Cells.Item(126, 14).Value = -10937  # N126

